$d = $word.ActiveDocument

# --- Change 1 ---
# "behavioral changes associated with shifting shelter-in-place policies have had"
# -> "behavioral changes associated with shelter-in-place restrictions and guidelines have had"
$d.Content.Find.Execute(
    "shifting shelter-in-place policies have had", $true, $false, $false, $false, $false,
    $true, 1, $false, "shelter-in-place restrictions and guidelines have had", 2) | Out-Null

# --- Change 2 ---
# " the response to extreme heat " -> " the mobility response to extreme heat "
$d.Content.Find.Execute(
    " the response to extreme heat ", $true, $false, $false, $false, $false,
    $true, 1, $false, " the mobility response to extreme heat ", 2) | Out-Null

# --- Change 3 ---
# "across the San Francisco Bay Area. We find" -> "across the region. We find"
$d.Content.Find.Execute(
    "across the San Francisco Bay Area. We find", $true, $false, $false, $false, $false,
    $true, 1, $false, "across the region. We find", 2) | Out-Null

# --- Change 4 ---
# "...indoor spaces have in heat mitigation and adaptation that allows for continued
#  activities in extreme temperatures. " with the _GoBack bookmark around "have "
# ->
# "...indoor spaces previously had in heat mitigation and adaptation that allows for
#  continued activities in extreme temperatures" + (empty _GoBack bookmark) + ". "
$d.Content.Find.Execute(
    "have in heat mitigation and adaptation that allows for continued activities in extreme temperatures. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "previously had in heat mitigation and adaptation that allows for continued activities in extreme temperatures. ",
    2) | Out-Null

# Reposition the _GoBack bookmark: it should now collapse to an empty range
# immediately before the final ". " at the end of the paragraph.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$found = $d.Content.Find.Execute("in extreme temperatures. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target = $d.Content.Find.Parent
    $periodPos = $target.End - 2
    $bmRange = $d.Range($periodPos, $periodPos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

Write-Host "done"
